# Small update to README
#
# The sheet had an AutoFilter on the "name" column (C) restricted to just
# "Binance USD", which hid every row whose name wasn't that. This change
# clears that filter criterion so all rows (7-66) are shown again, and
# updates the saved view state (zoom level + selected cell + window
# placement) to match what was active when the file was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the filter criteria on column 3 ("name") of the A1:J66 AutoFilter
# range. This removes the <filterColumn>/<filters> restriction and, as a
# consequence, un-hides every row that had been hidden by it.
$ws.Range("A1:J66").AutoFilter(3)

# Belt-and-braces: make sure nothing is left hidden.
$ws.ShowAllData()

# Update the window chrome: zoom level and the selected cell.
$excel.ActiveWindow.Zoom = 187
$ws.Range("C18").Select()

# Reposition / resize the workbook window to match the saved view state.
$win = $excel.ActiveWindow
$win.Left = 3120
$win.Top = 1700
$win.Width = 36300
$win.Height = 23240
